$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D7:D16").Value = "*"
$ws.Range("D19:D31").Value = "*"
